$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replace the checkmark/crossmark glyphs in the "Is Active" column ---
# F2:F5 originally read "✓ Active" (shared text), F6:F7 originally read "✗ Inactive".
# Write each group at once so the duplicated text collapses back to a single value.
$ws.Range("F2:F5").Value = "Active"
$ws.Range("F6:F7").Value = "Inactive"

# --- 2. Fix wrapping / alignment on the data rows (row 1 is the header, untouched) ---
# Columns A-E and G switch from "general" horizontal alignment to "left", and every
# data style turns wrapping on. Column F (the Active/Inactive flag) keeps its centered
# alignment but also gets wrapping turned on.
$ws.Range("A2:E7").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A2:E7").WrapText = $true

$ws.Range("G2:G7").HorizontalAlignment = -4131   # xlLeft
$ws.Range("G2:G7").WrapText = $true

$ws.Range("F2:F7").WrapText = $true

# --- 3. Widen the columns so the now-wrapped text still reads comfortably ---
# Every column grows by ~2 characters, except column D (Hire Date) which is untouched.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth + 2
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth + 2
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(3).ColumnWidth + 2
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(5).ColumnWidth + 2
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(6).ColumnWidth + 2
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(7).ColumnWidth + 2
